$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C9").Value = "Christian, Travis"
$ws.Range("C10").Value = "Christian, Travis"
$ws.Range("C11").Value = "Christian, Travis"

$ws.Range("H18").Select()
